# Update countries & provincias Spain
# Daily refresh of the COVID-19 country table: the "last updated" note,
# several countries' ranking order (swaps in column A), and the
# associated case/death/recovery figures (columns B-H) for the rows
# whose ranking moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 3 de Julio de 2020 a las 00:34'

$ws.Range('B4').Value = 2828738
$ws.Range('C4').Value = 48785
$ws.Range('D4').Value = 1185019
$ws.Range('E4').Value = 1512306
$ws.Range('G4').Value = 615
$ws.Range('H4').Value = 131413

$ws.Range('B5').Value = 1496858
$ws.Range('C5').Value = 43489
$ws.Range('E5').Value = 518827
$ws.Range('G5').Value = 1171
$ws.Range('H5').Value = 61884

$ws.Range('A8').Value = 'España'
$ws.Range('B8').Value = 297183
$ws.Range('C8').Value = 444
$ws.Range('G8').Value = 5
$ws.Range('H8').Value = 28368

$ws.Range('A9').Value = 'Peru'
$ws.Range('B9').Value = 292004
$ws.Range('C9').Value = 3527
$ws.Range('D9').Value = 182097
$ws.Range('E9').Value = 99862
$ws.Range('G9').Value = 185
$ws.Range('H9').Value = 10045

$ws.Range('A10').Value = 'Chile'
$ws.Range('B10').Value = 284541
$ws.Range('C10').Value = 2498
$ws.Range('D10').Value = 249247
$ws.Range('E10').Value = 29374
$ws.Range('G10').Value = 167
$ws.Range('H10').Value = 5920

$ws.Range('A11').Value = 'Reino Unido'
$ws.Range('B11').Value = 283757
$ws.Range('C11').Value = 576
$ws.Range('D11').Value = 0
$ws.Range('E11').Value = 0
$ws.Range('G11').Value = 89
$ws.Range('H11').Value = 43995

$ws.Range('B18').Value = 196706
$ws.Range('C18').Value = 382
$ws.Range('E18').Value = 7342
$ws.Range('G18').Value = 3
$ws.Range('H18').Value = 9064

$ws.Range('A19').Value = 'Sudafrica'
$ws.Range('B19').Value = 168061
$ws.Range('C19').Value = 8728
$ws.Range('D19').Value = 81999
$ws.Range('E19').Value = 83218
$ws.Range('G19').Value = 95
$ws.Range('H19').Value = 2844

$ws.Range('A20').Value = 'Francia'
$ws.Range('B20').Value = 166378
$ws.Range('C20').Value = 659
$ws.Range('D20').Value = 76802
$ws.Range('E20').Value = 59701
$ws.Range('G20').Value = 14
$ws.Range('H20').Value = 29875

$ws.Range('A22').Value = 'Colombia'
$ws.Range('B22').Value = 106110
$ws.Range('C22').Value = 4101
$ws.Range('D22').Value = 44531
$ws.Range('E22').Value = 57938
$ws.Range('G22').Value = 171
$ws.Range('H22').Value = 3641

$ws.Range('A23').Value = 'Canada'
$ws.Range('B23').Value = 104643
$ws.Range('C23').Value = 372
$ws.Range('D23').Value = 68217
$ws.Range('E23').Value = 27789
$ws.Range('G23').Value = 22
$ws.Range('H23').Value = 8637

$ws.Range('A27').Value = 'Argentina'
$ws.Range('B27').Value = 69941
$ws.Range('C27').Value = 2744
$ws.Range('D27').Value = 24186
$ws.Range('E27').Value = 44370
$ws.Range('G27').Value = 34
$ws.Range('H27').Value = 1385

$ws.Range('A28').Value = 'Suecia'
$ws.Range('B28').Value = 69692
$ws.Range('D28').Value = 0
$ws.Range('E28').Value = 0
$ws.Range('G28').Value = 0
$ws.Range('H28').Value = 5370

$ws.Range('A49').Value = 'Barein'
$ws.Range('B49').Value = 27837
$ws.Range('C49').Value = 423
$ws.Range('D49').Value = 22583
$ws.Range('E49').Value = 5160
$ws.Range('G49').Value = 2
$ws.Range('H49').Value = 94

$ws.Range('A50').Value = 'Rumania'
$ws.Range('B50').Value = 27746
$ws.Range('C50').Value = 450
$ws.Range('D50').Value = 19363
$ws.Range('E50').Value = 6696
$ws.Range('G50').Value = 20
$ws.Range('H50').Value = 1687

$ws.Range('B57').Value = 18874
$ws.Range('C57').Value = 151
$ws.Range('D57').Value = 16772
$ws.Range('E57').Value = 1127
$ws.Range('G57').Value = 1
$ws.Range('H57').Value = 975

$ws.Range('B72').Value = 9078
$ws.Range('C72').Value = 297
$ws.Range('D72').Value = 6034
$ws.Range('E72').Value = 3017

$ws.Range('B89').Value = 5315
$ws.Range('C89').Value = 161
$ws.Range('D89').Value = 2802
$ws.Range('E89').Value = 2281

$ws.Range('E100').Value = 2611
$ws.Range('G100').Value = 1
$ws.Range('H100').Value = 9

$ws.Range('B133').Value = 1081
$ws.Range('C133').Value = 6
$ws.Range('D133').Value = 959
$ws.Range('E133').Value = 54
$ws.Range('G133').Value = 1
$ws.Range('H133').Value = 68

$ws.Range('B154').Value = 547
$ws.Range('C154').Value = 12
$ws.Range('D154').Value = 252
$ws.Range('E154').Value = 282

$ws.Range('B166').Value = 250
$ws.Range('C166').Value = 2
$ws.Range('D166').Value = 117
$ws.Range('G166').Value = 1
$ws.Range('H166').Value = 14

$ws.Range('D171').Value = 194
$ws.Range('E171').Value = 6

$ws.Range('A203').Value = 'Laos'

$ws.Range('A204').Value = 'Santa Lucia'

$ws.Range('A205').Value = 'Dominica'

$ws.Range('A206').Value = 'Fiyi'

$ws.Range('A209').Value = 'Islas Malvinas'

$ws.Range('A210').Value = 'Groenlandia'
